$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 372, shifting existing rows (372..454) down to (373..455)
$ws.Rows("372:372").Insert()

# Populate the newly inserted row with the latest weekly price data
$ws.Range("A372").Value = 9
$ws.Range("B372").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C372").Value = "Metropolitana"
$ws.Range("D372").Value = 45209
$ws.Range("E372").Value = 13
$ws.Range("F372").Value = 100112043
$ws.Range("G372").Value = "Pepino ensalada"
$ws.Range("H372").Value = "Sin especificar"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 70
$ws.Range("K372").Value = 11000
$ws.Range("L372").Value = 12000
$ws.Range("M372").Value = 11500
$ws.Range("N372").Value = "$/caja 60 unidades"
$ws.Range("O372").Value = "Región de Arica y Parinacota"
$ws.Range("P372").Value = 192
$ws.Range("Q372").Value = 60
$ws.Range("R372").Value = "Hortaliza"
